$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 3.182878228561681
$ws.Range("C2").Value = 87981.0709163148
$ws.Range("D2").Value = 157.8057217802531
$ws.Range("E2").Value = 246.9852506941017
$ws.Range("G2").Value = 88389.04476701771
